# Switch to single faculty list in FacultyAvailabilityMatrix
# Adds a new row (row 14) to Sheet1 representing an additional faculty
# availability entry ("W/ngb/ngd"), mirroring the formatting used by the
# existing time-slot rows above it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row label (A14), formatted like A2:A13 but without the cell border ---
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A14").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("A14").Borders.LineStyle = -4142         # xlLineStyleNone
$ws.Range("A14").Value = "W/ngb/ngd"

# --- Plain (default-style) availability marks ---
$plainCols = @("M","O","Y","AH","AJ","AK","BC","BD","BM","BQ","BS","BU","BZ")
foreach ($col in $plainCols) {
    $ws.Range($col + "14").Value = 1
}

# --- Cells matching the borderless bold/centered numeric style (like AI3) ---
$ws.Range("AI3").Copy() | Out-Null
foreach ($col in @("U","AM","AN","BN","CA")) {
    $addr = $col + "14"
    $ws.Range($addr).PasteSpecial(-4122) | Out-Null  # xlPasteFormats
    $ws.Range($addr).Value = 1
}

# --- Cell matching column AS's own default style (inherited automatically) ---
$ws.Range("AS14").Value = 1

# --- Cell matching the borderless centered style (like U8) ---
$ws.Range("U8").Copy() | Out-Null
$ws.Range("BY14").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("BY14").Value = 1
